$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.890.59'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.741.13'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.98'
$ws.Range('E5').Value = '  -3.79%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5149'
$ws.Range('E7').Value = '  +1.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2824'
$ws.Range('E8').Value = '  +8.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.19'
$ws.Range('E9').Value = '  -3.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06093'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').Value = '1.743.24'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06983'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.38'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('E14').Value = '  +5.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.499'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.45'
$ws.Range('E16').Value = '  -2.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9992'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9990'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '25.887.52'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.51'
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006595'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '1.967.72'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.100'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.494'
$ws.Range('E24').Value = '  +3.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.125'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.68'
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.510'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.03'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.815'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.52'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08319'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.623'
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.408'
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04391'
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.619'
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9729'
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6073'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.669'
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01558'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.932'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9980'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.64'
$ws.Range('E42').Value = '  -2.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3840'
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7228'
$ws.Range('E44').Value = '  -3.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.936'
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05450'
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.357'
$ws.Range('E47').Value = '  +7.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1108'
$ws.Range('E48').Value = '  +3.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.55'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.80'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.521'
$ws.Range('E51').Value = '  +1.57%  '
